# AITATR3-265 Revisi financial nap4
# Update simulation input values on "Gross Yield (CF)" and "Regular Fixed" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: Gross Yield (CF) ---
$ws = $wb.Worksheets.Item("Gross Yield (CF)")

$ws.Range("B2").Value = 122200000          # was 3.5032E8
$ws.Range("B3").Value = 558116998.55       # was 558153000.49000001
$ws.Range("B5").Value = 610000000          # was 6.02E8
$ws.Range("B6").Value = 122500000          # was 3.5062E8
$ws.Range("E6").Value = 11134500           # was 4989500.0
$ws.Range("H6").Value = 10967346           # was 5022000.0
$ws.Range("E7").Value = 21440625           # was 0.0
$ws.Range("H7").Value = 22437500           # was 0.0
$ws.Range("B8").Value = 0.14682159         # was 0.14684329
$ws.Range("B9").Value = 12                 # was 24.0
$ws.Range("E15").Value = 0                 # was 0
$ws.Range("B16").Value = 0                 # was 0
$ws.Range("E16").Value = 3                 # was 3
$ws.Range("B17").Value = 62012999.84       # was 62017000.049999997
$ws.Range("D25").Value = 6432188           # was 6465615
$ws.Range("D26").Value = 100000            # was 100000
$ws.Range("D27").Value = 1000000           # was 1000000
$ws.Range("D28").Value = 150000            # was 150000
$ws.Range("D29").Value = 250000            # was 250000
$ws.Range("D31").Value = 40000             # was 40000
$ws.Range("D38").Value = 20000             # was 20000

# --- Sheet: Regular Fixed ---
$ws2 = $wb.Worksheets.Item("Regular Fixed")
$ws2.Range("B212").Value = 3               # was 3

Write-Host "Applied financial revision updates (AITATR3-265)."
